$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 2020-07-26 refresh of "Fonds de solidarite" VOLET1 regional/legal-category data.
# Each entry below is (row, nombre_aides, montant_total) reflecting the updated
# counts and amounts published for this extract; all other columns are unchanged.
$updates = @(
    @{ Row = 2; NombreAides = 306385; MontantTotal = 390931793 }
    @{ Row = 3; NombreAides = 253; MontantTotal = 301099 }
    @{ Row = 4; NombreAides = 305; MontantTotal = 436843 }
    @{ Row = 8; NombreAides = 822; MontantTotal = 1211544 }
    @{ Row = 10; NombreAides = 113104; MontantTotal = 165823559 }
    @{ Row = 12; NombreAides = 56673; MontantTotal = 81836194 }
    @{ Row = 16; NombreAides = 3755; MontantTotal = 5331173 }
    @{ Row = 20; NombreAides = 6006; MontantTotal = 8392260 }
    @{ Row = 22; NombreAides = 74514; MontantTotal = 93161583 }
    @{ Row = 28; NombreAides = 31511; MontantTotal = 46145854 }
    @{ Row = 30; NombreAides = 11044; MontantTotal = 15911020 }
    @{ Row = 33; NombreAides = 1478; MontantTotal = 2076871 }
    @{ Row = 35; NombreAides = 1648; MontantTotal = 2323974 }
    @{ Row = 36; NombreAides = 94012; MontantTotal = 118609288 }
    @{ Row = 37; NombreAides = 64; MontantTotal = 75477 }
    @{ Row = 38; NombreAides = 81; MontantTotal = 112488 }
    @{ Row = 44; NombreAides = 43264; MontantTotal = 63439607 }
    @{ Row = 46; NombreAides = 8802; MontantTotal = 12634147 }
    @{ Row = 48; NombreAides = 1336; MontantTotal = 1851857 }
    @{ Row = 51; NombreAides = 2103; MontantTotal = 2928539 }
    @{ Row = 52; NombreAides = 66325; MontantTotal = 83312356 }
    @{ Row = 53; NombreAides = 36; MontantTotal = 39883 }
    @{ Row = 54; NombreAides = 42; MontantTotal = 57866 }
    @{ Row = 58; NombreAides = 27317; MontantTotal = 40066703 }
    @{ Row = 59; NombreAides = 24; MontantTotal = 36000 }
    @{ Row = 61; NombreAides = 10623; MontantTotal = 15364049 }
    @{ Row = 63; NombreAides = 1309; MontantTotal = 1825706 }
    @{ Row = 67; NombreAides = 1358; MontantTotal = 1897181 }
    @{ Row = 69; NombreAides = 19541; MontantTotal = 25579770 }
    @{ Row = 73; NombreAides = 7174; MontantTotal = 10503921 }
    @{ Row = 75; NombreAides = 4843; MontantTotal = 7032802 }
    @{ Row = 76; NombreAides = 460; MontantTotal = 649689 }
    @{ Row = 77; NombreAides = 249; MontantTotal = 349904 }
    @{ Row = 78; NombreAides = 135230; MontantTotal = 168823155 }
    @{ Row = 79; NombreAides = 66; MontantTotal = 78110 }
    @{ Row = 82; NombreAides = 409; MontantTotal = 597475 }
    @{ Row = 84; NombreAides = 61513; MontantTotal = 90202951 }
    @{ Row = 87; NombreAides = 28480; MontantTotal = 41225419 }
    @{ Row = 89; NombreAides = 2553; MontantTotal = 3676471 }
    @{ Row = 90; NombreAides = 2541; MontantTotal = 3583312 }
    @{ Row = 91; NombreAides = 29757; MontantTotal = 40326921 }
    @{ Row = 95; NombreAides = 7451; MontantTotal = 10973440 }
    @{ Row = 97; NombreAides = 6706; MontantTotal = 9717049 }
    @{ Row = 99; NombreAides = 473; MontantTotal = 672805 }
    @{ Row = 100; NombreAides = 452; MontantTotal = 653139 }
    @{ Row = 101; NombreAides = 7605; MontantTotal = 10535403 }
    @{ Row = 103; NombreAides = 1932; MontantTotal = 2842255 }
    @{ Row = 105; NombreAides = 2621; MontantTotal = 3827867 }
    @{ Row = 107; NombreAides = 109; MontantTotal = 157620 }
    @{ Row = 108; NombreAides = 127; MontantTotal = 181689 }
    @{ Row = 109; NombreAides = 135891; MontantTotal = 168168386 }
    @{ Row = 113; NombreAides = 923; MontantTotal = 1355856 }
    @{ Row = 115; NombreAides = 51270; MontantTotal = 75187817 }
    @{ Row = 117; NombreAides = 25916; MontantTotal = 37553759 }
    @{ Row = 118; NombreAides = 1233; MontantTotal = 1685352 }
    @{ Row = 121; NombreAides = 2055; MontantTotal = 2887981 }
    @{ Row = 123; NombreAides = 468356; MontantTotal = 617507011 }
    @{ Row = 124; NombreAides = 86; MontantTotal = 114117 }
    @{ Row = 128; NombreAides = 1306; MontantTotal = 1937238 }
    @{ Row = 130; NombreAides = 198029; MontantTotal = 291240890 }
    @{ Row = 131; NombreAides = 366; MontantTotal = 545790 }
    @{ Row = 133; NombreAides = 169855; MontantTotal = 246993802 }
    @{ Row = 136; NombreAides = 2648; MontantTotal = 3718731 }
    @{ Row = 138; NombreAides = 5652; MontantTotal = 7980217 }
    @{ Row = 141; NombreAides = 42126; MontantTotal = 56330028 }
    @{ Row = 145; NombreAides = 5; MontantTotal = 7500 }
    @{ Row = 147; NombreAides = 13530; MontantTotal = 19860035 }
    @{ Row = 148; NombreAides = 3564; MontantTotal = 5141043 }
    @{ Row = 151; NombreAides = 370; MontantTotal = 532222 }
    @{ Row = 153; NombreAides = 347; MontantTotal = 487809 }
    @{ Row = 154; NombreAides = 16412; MontantTotal = 21742798 }
    @{ Row = 158; NombreAides = 6791; MontantTotal = 9885802 }
    @{ Row = 160; NombreAides = 4625; MontantTotal = 6657989 }
    @{ Row = 163; NombreAides = 248; MontantTotal = 355433 }
    @{ Row = 165; NombreAides = 13659; MontantTotal = 19817831 }
    @{ Row = 166; NombreAides = 1633; MontantTotal = 2429138 }
    @{ Row = 167; NombreAides = 215; MontantTotal = 317302 }
    @{ Row = 170; NombreAides = 79; MontantTotal = 118449 }
    @{ Row = 171; NombreAides = 84607; MontantTotal = 105900308 }
    @{ Row = 176; NombreAides = 627; MontantTotal = 924346 }
    @{ Row = 178; NombreAides = 32890; MontantTotal = 48249164 }
    @{ Row = 180; NombreAides = 12475; MontantTotal = 18035889 }
    @{ Row = 182; NombreAides = 1181; MontantTotal = 1652676 }
    @{ Row = 184; NombreAides = 1497; MontantTotal = 2106283 }
    @{ Row = 186; NombreAides = 228722; MontantTotal = 284573779 }
    @{ Row = 194; NombreAides = 84161; MontantTotal = 123393817 }
    @{ Row = 197; NombreAides = 31726; MontantTotal = 45670091 }
    @{ Row = 200; NombreAides = 4793; MontantTotal = 6827260 }
    @{ Row = 203; NombreAides = 4413; MontantTotal = 6121390 }
    @{ Row = 205; NombreAides = 11; MontantTotal = 13411 }
    @{ Row = 206; NombreAides = 251765; MontantTotal = 311746683 }
    @{ Row = 213; NombreAides = 601; MontantTotal = 875856 }
    @{ Row = 215; NombreAides = 91760; MontantTotal = 134293926 }
    @{ Row = 218; NombreAides = 49134; MontantTotal = 71060029 }
    @{ Row = 219; NombreAides = 29; MontantTotal = 40769 }
    @{ Row = 221; NombreAides = 4424; MontantTotal = 6206721 }
    @{ Row = 224; NombreAides = 5206; MontantTotal = 7194064 }
    @{ Row = 227; NombreAides = 102198; MontantTotal = 128113851 }
    @{ Row = 232; NombreAides = 553; MontantTotal = 808839 }
    @{ Row = 234; NombreAides = 48151; MontantTotal = 70568176 }
    @{ Row = 236; NombreAides = 11878; MontantTotal = 17079547 }
    @{ Row = 238; NombreAides = 1812; MontantTotal = 2600513 }
    @{ Row = 240; NombreAides = 2335; MontantTotal = 3263033 }
    @{ Row = 241; NombreAides = 246297; MontantTotal = 311307350 }
    @{ Row = 242; NombreAides = 162; MontantTotal = 201831 }
    @{ Row = 243; NombreAides = 240; MontantTotal = 344457 }
    @{ Row = 249; NombreAides = 92580; MontantTotal = 135705975 }
    @{ Row = 252; NombreAides = 62106; MontantTotal = 90050353 }
    @{ Row = 254; NombreAides = 2305; MontantTotal = 3254042 }
    @{ Row = 257; NombreAides = 4186; MontantTotal = 5874601 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.NombreAides
    $ws.Cells.Item($u.Row, 4).Value = $u.MontantTotal
}
